$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: D3 had the 'MEC-1B-T. M. Metalicos' entry, move it to F3; D3 becomes "-"
$ws.Cells.Item(3, 4).Value2 = "-"
$ws.Cells.Item(3, 5).Value2 = "-"
$ws.Cells.Item(3, 6).Value2 = "['MEC-1B-T. M. Metalicos', 'MEC-1B-T. M. Metalicos']"

# Row 11: B11 gets the 'MEC-1A-T. M. Metalicos' entry; the Metalografia list moves from E11 to F11
$ws.Cells.Item(11, 2).Value2 = "['MEC-1A-T. M. Metalicos', 'MEC-1A-T. M. Metalicos']"
$ws.Cells.Item(11, 5).Value2 = "-"
$ws.Cells.Item(11, 6).Value2 = "[-, 'MEC-2A-Metalografia', -, -]"

# Row 12: Metalografia list moves from E12 to F12 (replacing the MEC-1A entry there)
$ws.Cells.Item(12, 5).Value2 = "-"
$ws.Cells.Item(12, 6).Value2 = "[-, 'MEC-2A-Metalografia', -, -]"

# Row 14: Metalografia list moves from E14 to F14
$ws.Cells.Item(14, 5).Value2 = "-"
$ws.Cells.Item(14, 6).Value2 = "[-, 'MEC-2A-Metalografia', -, -]"

# Row 15: Metalografia list moves from E15 to F15
$ws.Cells.Item(15, 5).Value2 = "-"
$ws.Cells.Item(15, 6).Value2 = "[-, 'MEC-2A-Metalografia', -, -]"
